# Project 4 / Task 3: split the old "checkInventory" row into a new
# "addInventory2" negative-value test case (re-using the row that used to
# hold checkInventory) and a fresh row right below it that keeps the
# original checkInventory content (pushed down, now starting a new
# rendered page).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wdns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Locate the existing "checkInventory" row.
# ---------------------------------------------------------------------
$checkRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $r = $t.Rows.Item($i)
    if ($r.Cells.Item(1).Range.Text -like "checkInventory*") {
        $checkRow = $r
        break
    }
}

$rowIndex = $checkRow.Index

# ---------------------------------------------------------------------
# 1) Insert a brand-new blank row right after it; Word clones the tcPr
#    (widths/borders/vAlign/hideMark) of the surrounding rows.
# ---------------------------------------------------------------------
$belowRow = $t.Rows.Item($rowIndex + 1)
$newRow = $t.Rows.Add($belowRow)

# ---------------------------------------------------------------------
# 2) Populate the new row with the ORIGINAL checkInventory content,
#    adding the lastRenderedPageBreak that now precedes it.
# ---------------------------------------------------------------------
$nc1 = $newRow.Cells.Item(1)
$nc1.Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>checkInventory</w:t></w:r></w:p>
"@) | Out-Null

$nc2 = $newRow.Cells.Item(2)
$nc2.Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Precondition: addInventory1 has run successfully.</w:t></w:r></w:p><w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Enter: Menu option </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>5,`u{201d}Check inventory`u{201d}</w:t></w:r></w:p><w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Return to main menu.</w:t></w:r></w:p>
"@) | Out-Null

$nc3 = $newRow.Cells.Item(3)
$nc3.Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="24"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>The current inventory is displayed successfully.</w:t></w:r></w:p><w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="24"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Coffee: 16, Milk: 17, Sugar: 16, Chocolate: 17</w:t></w:r></w:p>
"@) | Out-Null

# 4th (blank, smaller font) cell already matches - leave as-is.

# ---------------------------------------------------------------------
# 3) Rewrite the original row in place to become "addInventory2".
# ---------------------------------------------------------------------

# give it the new trHeight (val=65 twips => 3.25 pt, Word stores pt*20)
$checkRow.Height = 3.25

# Cell 1: rename.
$checkRow.Cells.Item(1).Range.Find.Execute("checkInventory", $true, $false, $false, $false, $false, $true, 1, $false, "addInventory2", 2) | Out-Null

# Cell 2: description - needs a 4th paragraph, then rewrite all four.
$c2 = $checkRow.Cells.Item(2)
$c2.Range.Paragraphs.Item($c2.Range.Paragraphs.Count).Range.InsertParagraphAfter() | Out-Null

$c2.Range.Paragraphs.Item(1).Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Precondition: Run CoffeeMaker</w:t></w:r></w:p>
"@) | Out-Null

$c2.Range.Paragraphs.Item(2).Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Enter: Menu option 4, `u{201c}Add inventory`u{201d}</w:t></w:r></w:p>
"@) | Out-Null

$c2.Range.Paragraphs.Item(3).Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Enter Coffee Amount: -1</w:t></w:r></w:p>
"@) | Out-Null

$c2.Range.Paragraphs.Item(4).Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>Return to the main menu.</w:t></w:r></w:p>
"@) | Out-Null

# Cell 3: expected results - drop down to a single paragraph, then rewrite it.
$c3 = $checkRow.Cells.Item(3)
$c3.Range.Paragraphs.Item(2).Range.Delete() | Out-Null
$c3.Range.Paragraphs.Item(1).Range.InsertXML(@"
<w:p $wdns><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="24"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t>Inventory is not added.</w:t></w:r></w:p>
"@) | Out-Null

# Cell 4: blank placeholder - unchanged.

Write-Output "done"
